$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '89.695.56'
$r.Style = "Normal"
$ws.Range("E2").Value = '  +1.78%  '
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '3.182.33'
$r.Style = "Normal"
$ws.Range("E3").Value = '  -2.57%  '
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '213.38'
$r.Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '614.89'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -2.51%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.391'
$r.Style = "Normal"
$ws.Range("E7").Value = '  +2.77%  '
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.690'
$r.Style = "Normal"
$ws.Range("E8").Value = '  -5.29%  '
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range("E9").Value = '  -0.03%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '3.176.44'
$r.Style = "Normal"
$ws.Range("E10").Value = '  -2.63%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.577'
$r.Style = "Normal"
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("E12").Value = '  -5.25%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.0000254'
$r.Style = "Normal"
$ws.Range("E13").Value = '  -3.85%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '89.509.51'
$r.Style = "Normal"
$ws.Range("E14").Value = '  +1.98%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '3.764.64'
$r.Style = "Normal"
$ws.Range("E15").Value = '  -2.66%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '32.90'
$r.Style = "Normal"
$ws.Range("E16").Value = '  -4.58%  '
$ws.Range("B17").Value = 'Toncoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '5.24'
$r.Style = "Normal"
$ws.Range("E17").Value = '  -5.29%  '
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '3.163.64'
$r.Style = "Normal"
$ws.Range("E18").Value = '  -4.14%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '3.28'
$r.Style = "Normal"
$ws.Range("E19").Value = '  +3.92%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '13.43'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -4.85%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '437.06'
$r.Style = "Normal"
$ws.Range("E21").Value = '  -0.29%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '0.0000195'
$r.Style = "Normal"
$ws.Range("E22").Value = '  +38.24%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '8.61'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -4.07%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '5.06'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -5.62%  '
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '5.14'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -2.39%  '
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '11.70'
$r.Style = "Normal"
$ws.Range("E26").Value = '  -5.16%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '3.336.46'
$r.Style = "Normal"
$ws.Range("E27").Value = '  -3.33%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '75.24'
$r.Style = "Normal"
$ws.Range("E28").Value = '  -2.95%  '
$ws.Range("E29").Value = '  +0.22%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '0.167'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -7.05%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E31").Value = '  -0.03%  '
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '4.18'
$r.Style = "Normal"
$ws.Range("E32").Value = '  +35.01%  '
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '8.44'
$r.Style = "Normal"
$ws.Range("E33").Value = '  -4.71%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '533.01'
$r.Style = "Normal"
$ws.Range("E34").Value = '  -6.16%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '7.01'
$r.Style = "Normal"
$ws.Range("E35").Value = '  -2.34%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '1.86'
$r.Style = "Normal"
$ws.Range("E36").Value = '  -5.85%  '
$ws.Range("E37").Value = '  -7.71%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '21.95'
$r.Style = "Normal"
$ws.Range("E38").Value = '  -4.38%  '
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '22.30'
$r.Style = "Normal"
$ws.Range("E39").Value = '  +2.21%  '
$ws.Range("E40").Value = '  -8.87%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  -0.04%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '1.92'
$r.Style = "Normal"
$ws.Range("E43").Value = '  -5.48%  '
$ws.Range("E44").Value = '  -7.39%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '150.23'
$r.Style = "Normal"
$ws.Range("E45").Value = '  -1.53%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '172.54'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -3.89%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '43.47'
$r.Style = "Normal"
$ws.Range("E47").Value = '  -2.79%  '
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '0.124'
$r.Style = "Normal"
$ws.Range("E48").Value = '  -9.58%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '1.24'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -7.06%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '4.07'
$r.Style = "Normal"
$ws.Range("E50").Value = '  -4.23%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '0.610'
$r.Style = "Normal"
$ws.Range("E51").Value = '  -3.12%  '

Write-Host "Update complete"